# RUN ID 1648 - fixed retrieval of config info from AH / fixed uploading of
# info to AH. The last data row (row 12) was a stray duplicate of row 11;
# replace it with the correct, new record that should have been uploaded.
#
# NOTE: every column in this sheet stores its data as shared-string TEXT
# (even the numeric-looking id/amount columns and the date column), never
# as a native Excel number/date. Plain `Range.Value = "830422"` would get
# auto-coerced by Excel into a real number (and a date-like string into a
# serial date), which would not match the source data. To force text entry
# without Excel's auto-detection kicking in, we prefix the literal with a
# leading apostrophe (the normal Excel "treat as text" escape) and then
# strip the resulting quote-prefix cell format so the cell's style is left
# exactly as it was (General, default style) - only the stored value/type
# changes, same as the original edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "'830422"
$ws.Range("B12").Value = "Professional Services"
$ws.Range("C12").Value = "'252934"
$ws.Range("D12").Value = "'50586.8"
$ws.Range("E12").Value = "'303521"
$ws.Range("F12").Value = "CAD"
$ws.Range("G12").Value = "'2017-12-19"

# Drop the quote-prefix formatting that typing an apostrophe adds, so the
# cells fall back to the same (default) style they had before the edit.
$ws.Range("A12:G12").ClearFormats()
